# Weekly update: insert a new "Zanahoria" (Femacal de La Calera) price
# observation as a new row 395 in the data table, pushing the existing
# rows 395:423 down to 396:424 (dimension grows from A1:R423 to A1:R424).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 395; Excel shifts rows 395:423 down to 396:424
# and copies the formatting (incl. the date style on column D) from the
# row above, same as a manual "Insert Row" in the UI.
$ws.Rows.Item(395).Insert()

# Populate the newly inserted row with this week's observation.
$ws.Cells.Item(395, 1).Value = 3
$ws.Cells.Item(395, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(395, 3).Value = "Coquimbo"
$ws.Cells.Item(395, 4).Value = 44826
$ws.Cells.Item(395, 5).Value = 5
$ws.Cells.Item(395, 6).Value = 100114013
$ws.Cells.Item(395, 7).Value = "Zanahoria"
$ws.Cells.Item(395, 8).Value = "Sin especificar"
$ws.Cells.Item(395, 9).Value = "Primera"
$ws.Cells.Item(395, 10).Value = 120
$ws.Cells.Item(395, 11).Value = 13000
$ws.Cells.Item(395, 12).Value = 13000
$ws.Cells.Item(395, 13).Value = 13000
$ws.Cells.Item(395, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(395, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(395, 16).Value = 650
$ws.Cells.Item(395, 17).Value = 20
$ws.Cells.Item(395, 18).Value = "Hortaliza"
